# The "象の群れ。ドローンの音から逃げている" post (row 810) was removed from
# the source data. Delete that entire worksheet row; Excel will shift every
# subsequent row up by one and update the sheet's used-range dimension
# automatically (from A1:C887 to A1:C886).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(810).Delete()
